$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "dammu" row (old row 3) entirely - this shifts the RRR row
# up from row 4 to row 3.
$ws.Rows(3).Delete()

# Conflict resolution in admin ratings (Admin Rating / column M).
$ws.Range("M2").Value = 38
$ws.Range("M3").Value = 97

# Restore the view to the selection left after resolving the conflict.
$ws.Range("F5").Select()
